# ESCALETA_CN_09_01_CO.xlsx — "Inclusion de proyectos metodo cientifico"
#
# A new row describing a "Competencias" project resource is inserted right
# before the current row 34 ("Fin de unidad" / Mapa conceptual). All the
# rows that used to occupy 34-36 move down to 35-37, and the trailing blank
# row moves from 37 to 38.
#
# The per-column cell styles are already identical across this whole block
# of rows (34-36 before the edit), so instead of using a native row-insert
# (which, on this host, stamps freshly allocated style ids on the inserted
# row instead of reusing the existing per-column ones) we copy the existing
# formatting down with Copy/PasteSpecial and then overwrite the cell
# contents directly — this keeps every cell's style id exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESCALETA")

# ---------------------------------------------------------------------
# 1. Make room: row 37 currently only holds a lone formatted-but-empty
#    "O37" cell (the sheet's trailing blank row). Give it the same
#    per-column formatting as row 36 has, then do the same for the new
#    row 38 (copying row 37's original formatting, i.e. just style 10 on
#    column O) before we overwrite any values.
# ---------------------------------------------------------------------
$ws.Range("O37").Copy()
$ws.Range("O38").PasteSpecial(-4122)

$ws.Range("A36:U36").Copy()
$ws.Range("A37:U37").PasteSpecial(-4122)

$ws.Rows.Item(37).RowHeight = $ws.Rows.Item(36).RowHeight()
$ws.Rows.Item(36).RowHeight = $ws.Rows.Item(35).RowHeight()

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Shift the old content of rows 34-36 down to rows 35-37 (values only;
#    formatting is already correct since these rows all shared the same
#    per-column styles before the edit).
# ---------------------------------------------------------------------

# old row 36 -> row 37
$ws.Range("D37").Value = "Fin de unidad"
$ws.Range("G37").Value = "Banco de actividades: El material hereditario y su expresión"
$ws.Range("H37").Value = 35
$ws.Range("J37").Value = "Motor que incluye preguntas de respuesta abierta del tema El material hereditario y su expresión"
$ws.Range("L37").Value = "ACTIVIDAD "
$ws.Range("M37").Value = ""
$ws.Range("N37").Value = "M101AP"
$ws.Range("O37").Value = ""
$ws.Range("P37").Value = "NO"
$ws.Range("Q37").Value = 6
$ws.Range("R37").Value = "RM"
$ws.Range("S37").Value = "Recursos M"
$ws.Range("T37").Value = "Recurso M101AP-01"
$ws.Range("U37").Value = "RM_01_01_CO"

# old row 35 -> row 36
$ws.Range("D36").Value = "Fin de unidad"
$ws.Range("G36").Value = "Evaluación"
$ws.Range("H36").Value = 34
$ws.Range("J36").Value = "Evalúa tus conocimientos acerca del tema El material hereditario y su expresión"
$ws.Range("L36").Value = "ACTIVIDAD "
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = "M4A"
$ws.Range("O36").Value = "Preguntas generales sobre lo visto en el capítulo"
$ws.Range("P36").Value = "SI"
$ws.Range("Q36").Value = 6
$ws.Range("R36").Value = "RM"
$ws.Range("S36").Value = "Recursos M"
$ws.Range("T36").Value = "Recurso M4A-03"
$ws.Range("U36").Value = "RM_01_01_CO"

# old row 34 -> row 35
$ws.Range("D35").Value = "Fin de unidad"
$ws.Range("G35").Value = "Mapa conceptual"
$ws.Range("H35").Value = 33
$ws.Range("J35").Value = "Mapa conceptual del tema El material hereditario y su expresión"
$ws.Range("L35").Value = ""
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = ""
$ws.Range("O35").Value = ""
$ws.Range("P35").Value = "SI"
$ws.Range("Q35").Value = ""
$ws.Range("R35").Value = ""
$ws.Range("S35").Value = ""
$ws.Range("T35").Value = ""
$ws.Range("U35").Value = ""

# ---------------------------------------------------------------------
# 3. Write the brand-new row 34: the "Competencias" project resource.
# ---------------------------------------------------------------------
$ws.Range("A34").Value = "CIENCIAS"
$ws.Range("B34").Value = "CN_09_01_CO"
$ws.Range("C34").Value = "El material hereditario y su expresión"
$ws.Range("D34").Value = "Competencias"
$ws.Range("G34").Value = "Competencias: analizar las variables de un proyecto de investigación "
$ws.Range("H34").Value = 32
$ws.Range("I34").Value = "No"
$ws.Range("J34").Value = "Proyecto que permite adquirir destrezas en el análisis de variables a partir de la utilización del método científico"
$ws.Range("K34").Value = "NUEVO "
$ws.Range("L34").Value = "INTERACTIVO "
$ws.Range("M34").Value = "F13"
$ws.Range("P34").Value = "NO"
$ws.Range("Q34").Value = 6
$ws.Range("R34").Value = "RF"
$ws.Range("S34").Value = "Recursos F"
$ws.Range("T34").Value = "Recurso F13-04"
$ws.Range("U34").Value = "RF_01_01_CO"

Write-Host "Done applying Competencias row insertion"
